$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C3").Value = "Shoulder press"

$ws.Range("C4").Value = "Tricep press"
$ws.Range("D4").Value = 70
$ws.Range("E4").Value = ""

$ws.Range("B5").Value = "x"
$ws.Range("C5").Value = "Lateral raises"
$ws.Range("D5").Value = 10
$ws.Range("E5").Value = ""

$ws.Range("A2").Value = "Mon"
$ws.Range("A7").Value = "Wed"
$ws.Range("A12").Value = "Fri"

$ws.Range("A13").Select()
